# Insert a new weekly record for "Acelga" at row 614, pushing the existing
# rows (614-661) down by one (to 615-662). This mirrors the new dimension
# A1:R662 from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(614).Insert()

$ws.Range("A614").Value = 3
$ws.Range("B614").Value = "Femacal de La Calera"
$ws.Range("C614").Value = "Coquimbo"
$ws.Range("D614").Value = 45265
$ws.Range("E614").Value = 5
$ws.Range("F614").Value = 100112009
$ws.Range("G614").Value = "Acelga"
$ws.Range("H614").Value = "Sin especificar"
$ws.Range("I614").Value = "Primera"
$ws.Range("J614").Value = 210
$ws.Range("K614").Value = 3000
$ws.Range("L614").Value = 3500
$ws.Range("M614").Value = 3262
$ws.Range("N614").Value = '$/docena de atados (6 kilos)'
$ws.Range("O614").Value = "Provincia de Quillota"
$ws.Range("P614").Value = 544
$ws.Range("Q614").Value = 6
$ws.Range("R614").Value = "Hortaliza"
